# Regenerate save_data to use K (kurtosis?) instead of Strike# in column G,
# writing the newly computed s_vals for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new value for column G ("K")
$updates = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 0
    6  = 2
    7  = 3
    9  = 2
    10 = 2
    12 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
